$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix G30: "Preposicao" -> "Substantivo"
$ws.Range("G30").Value = 'Substantivo'

# Add new question rows 32-40 (ids 31-39)
# Row 32
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 'Cesgranrio'
$ws.Range("C32").Value = 'UNIRIO'
$ws.Range("D32").Value = '''2016'
$ws.Range("E32").Value = 'Considere-se a seguinte passagem de Texto: “Dir-se-ia que os blocos <u>impostos </u>são opressivos e obrigatórios”
A classe da palavra <u>impostos </u>no trecho acima é a mesma da palavra destacada em:'
$ws.Range("F32").Value = 'Português'
$ws.Range("G32").Value = 'Substantivo'
$ws.Range("H32").Value = 'Médio'
$ws.Range("I32").Value = 'ME'
$ws.Range("J32").Value = 'O Congresso debateu muito, mas autorizou o aumento do <u>imposto </u>de renda.'
$ws.Range("K32").Value = 'Muitas pessoas se impressionam com qualquer estilo <u>imposto </u>pela mídia.'
$ws.Range("L32").Value = 'A enfermeira chegou logo a seguir de um grito esganiçado que foi <u>imposto</u> pelo futuro pai.'
$ws.Range("M32").Value = 'A mudança da moda é o <u>imposto </u>que a indústria do pobre lança sobre a vaidade do rico.'
$ws.Range("N32").Value = 'O padre tinha <u>imposto </u>uma pesada penitência àquele infeliz pecador.'
$ws.Range("O32").Value = 'B'
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 'Cesgranrio'
$ws.Range("C33").Value = 'Liquigás'
$ws.Range("D33").Value = '''2014'
$ws.Range("E33").Value = 'Ocorre a formação do plural de maneira idêntica à que acontece com a palavra irmão em'
$ws.Range("F33").Value = 'Português'
$ws.Range("G33").Value = 'Substantivo'
$ws.Range("H33").Value = 'Fácil'
$ws.Range("I33").Value = 'ME'
$ws.Range("J33").Value = 'aproximação'
$ws.Range("K33").Value = 'alemão'
$ws.Range("L33").Value = 'cirurgião'
$ws.Range("M33").Value = 'órgão'
$ws.Range("N33").Value = 'guardião'
$ws.Range("O33").Value = 'D'
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 'Cesgranrio'
$ws.Range("C34").Value = 'Banco do Brasil'
$ws.Range("D34").Value = '''2014'
$ws.Range("E34").Value = 'No fragmento “fazer um safári, frequentar uma praia de nudismo, comer algo exótico (um baiacu venenoso, por exemplo), visitar um vulcão ativo”, são palavras de classes gramaticais diferentes'
$ws.Range("F34").Value = 'Português'
$ws.Range("G34").Value = 'Substantivo'
$ws.Range("H34").Value = 'Fácil'
$ws.Range("I34").Value = 'ME'
$ws.Range("J34").Value = '“praia” e “ativo”'
$ws.Range("K34").Value = '“venenoso” e “exótico"'
$ws.Range("L34").Value = '“baiacu” e “nudismo”'
$ws.Range("M34").Value = '“ativo” e “exótico"'
$ws.Range("N34").Value = '“safári” e “vulcão”'
$ws.Range("O34").Value = 'A'
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 'Cesgranrio'
$ws.Range("C35").Value = 'Liquigás'
$ws.Range("D35").Value = '''2013'
$ws.Range("E35").Value = 'No trecho do Texto “capazes de sentar, interagir e <u>celebrar</u> com nossos semelhantes.”, o verbo destacado dá origem ao substantivo derivado celebração, grafado com ç. Os dois verbos que formam substantivos derivados grafados com ç são'
$ws.Range("F35").Value = 'Português'
$ws.Range("G35").Value = 'Substantivo'
$ws.Range("H35").Value = 'Médio'
$ws.Range("I35").Value = 'ME'
$ws.Range("J35").Value = 'combinar, nomear'
$ws.Range("K35").Value = 'elaborar, agredir'
$ws.Range("L35").Value = 'permitir, denominar'
$ws.Range("M35").Value = 'progredir, coroar'
$ws.Range("N35").Value = 'trair, compreender'
$ws.Range("O35").Value = 'A'
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 'Cesgranrio'
$ws.Range("C36").Value = 'BNDES'
$ws.Range("D36").Value = '''2013'
$ws.Range("E36").Value = 'Faz o plural como palavra-chave, com dupla possibilidade de flexão, o composto'
$ws.Range("F36").Value = 'Português'
$ws.Range("G36").Value = 'Substantivo'
$ws.Range("H36").Value = 'Difícil'
$ws.Range("I36").Value = 'ME'
$ws.Range("J36").Value = 'lugar-comum'
$ws.Range("K36").Value = 'guarda-roupa'
$ws.Range("L36").Value = 'aço-liga'
$ws.Range("M36").Value = 'amor-perfeito'
$ws.Range("N36").Value = 'abaixo-assinado'
$ws.Range("O36").Value = 'C'
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0

# Row 37
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 'Cesgranrio'
$ws.Range("C37").Value = 'Petrobras'
$ws.Range("D37").Value = '''2014'
$ws.Range("E37").Value = 'O fragmento do texto em que o vocábulo em destaque foi substantivado é:'
$ws.Range("F37").Value = 'Português'
$ws.Range("G37").Value = 'Substantivo'
$ws.Range("H37").Value = 'Difícil'
$ws.Range("I37").Value = 'ME'
$ws.Range("J37").Value = '“sua imagem foi literalmente apagada de fotografias dos líderes da revolução”.'
$ws.Range("K37").Value = '“A técnica usada para eliminar o Trotsky”.'
$ws.Range("L37").Value = '“Existe até uma técnica para retocar a imagem em movimento”.'
$ws.Range("M37").Value = '“Se a prova fotográfica não vale mais nada nestes novos tempos inconfiáveis, a assinatura muito menos”.'
$ws.Range("N37").Value = '“E se eu estiver fazendo a barba e escovando os dentes de um impostor, de um eu apócrifo?”.'
$ws.Range("O37").Value = 'A'
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0

# Row 38
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 'Cesgranrio'
$ws.Range("C38").Value = 'Liquigás'
$ws.Range("D38").Value = '''2018'
$ws.Range("E38").Value = 'Observe a expressão “velhas casas brasileiras”.
Caso o redator tivesse escrito “casas velhas brasileiras”, o trecho'
$ws.Range("F38").Value = 'Português'
$ws.Range("G38").Value = 'Adjetivo'
$ws.Range("H38").Value = 'Médio'
$ws.Range("I38").Value = 'ME'
$ws.Range("J38").Value = 'permaneceria com o mesmo sentido.'
$ws.Range("K38").Value = 'indicaria que as casas estavam abandonadas.'
$ws.Range("L38").Value = 'mostraria as casas como construções populares.'
$ws.Range("M38").Value = 'inverteria o sentido de <u>casas </u>e de <u>velhas</u>.'
$ws.Range("N38").Value = 'passaria a indicar as casas como gastas pelo tempo.'
$ws.Range("O38").Value = 'E'
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 0

# Row 39
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 'Cesgranrio'
$ws.Range("C39").Value = 'CEFET-RJ'
$ws.Range("D39").Value = '''2014'
$ws.Range("E39").Value = 'No trecho do Texto “<u>Só </u>de experimentar fiquei suando”, a palavra destacada tem o mesmo sentido em:'
$ws.Range("F39").Value = 'Português'
$ws.Range("G39").Value = 'Adjetivo'
$ws.Range("H39").Value = 'Fácil'
$ws.Range("I39").Value = 'ME'
$ws.Range("J39").Value = 'Não o deixem <u>só</u>!'
$ws.Range("K39").Value = 'Ele andou <u>só</u>.'
$ws.Range("L39").Value = 'Eles compraram <u>só </u>pelo prazer de comprar.'
$ws.Range("M39").Value = 'ser humano está muito <u>só</u>.'
$ws.Range("N39").Value = 'O homem sentiu-se <u>só </u>naquele lugar frio.'
$ws.Range("O39").Value = 'C'
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0

# Row 40
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 'Cesgranrio'
$ws.Range("C40").Value = 'CEFET-RJ'
$ws.Range("D40").Value = '''2014'
$ws.Range("E40").Value = 'Em qual dos períodos abaixo, a troca de posição entre a palavra sublinhada e o substantivo a que se refere mantém o sentido?'
$ws.Range("F40").Value = 'Português'
$ws.Range("G40").Value = 'Adjetivo'
$ws.Range("H40").Value = 'Fácil'
$ws.Range("I40").Value = 'ME'
$ws.Range("J40").Value = '<u>Algum </u>autor desejava a minha opinião sobre o seu trabalho.'
$ws.Range("K40").Value = 'O <u>mesmo </u>porteiro me entregou o pacote na recepção do hotel.'
$ws.Range("L40").Value = 'Meu pai procurou uma <u>certa </u>pessoa para me entregar o embrulho.'
$ws.Range("M40").Value = 'Contar histórias é uma <u>prazerosa </u>forma de aproximar os indivíduos.'
$ws.Range("N40").Value = '<u>Grandes </u>poemas épicos servem para perpetuar a cultura de um povo.'
$ws.Range("O40").Value = 'D'
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 0
